$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 142, shifting existing rows 142:174 down to 143:175
$ws.Rows(142).Insert()

# Populate the newly inserted row 142 with the new record.
$ws.Range("A142").Value = 4
$ws.Range("B142").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C142").Value = "Los Lagos"
$ws.Range("D142").Value = "2021-11-11"
$ws.Range("E142").Value = 10
$ws.Range("F142").Value = 100112043
$ws.Range("G142").Value = "Pepino ensalada"
$ws.Range("H142").Value = "Sin especificar"
$ws.Range("I142").Value = "Primera"
$ws.Range("J142").Value = 200
$ws.Range("K142").Value = 12000
$ws.Range("L142").Value = 12000
$ws.Range("M142").Value = 12000
$ws.Range("N142").Value = "`$/caja 60 unidades"
$ws.Range("O142").Value = "Región de Arica y Parinacota"
$ws.Range("P142").Value = 200
$ws.Range("Q142").Value = 60
$ws.Range("R142").Value = "Hortaliza"
